$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-6 (Date, Match, Toss Winner, Match Winner)
$data = @(
    @("23-03-2025", "Sunrisers Hyderabad vs Rajasthan Royals", "Rajasthan Royals", "Rajasthan Royals"),
    @("25-03-2025", "Chennai Super Kings vs Mumbai Indians", "Mumbai Indians", "Mumbai Indians"),
    @("25-03-2025", "Gujarat Titans vs Punjab Kings", "Gujarat Titans", "Gujarat Titans"),
    @("24-03-2025", "Delhi Capitals vs Lucknow Super Giants", "Delhi Capitals", "Delhi Capitals"),
    @("22-03-2025", "Kolkata Knight Riders vs Royal Challengers Bengaluru", "Royal Challengers Bengaluru", "Kolkata Knight Riders")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
